$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.064.77'
$ws.Range("E2").Value = '  -3.66%  '
$ws.Range("D3").Value = '1.643.19'
$ws.Range("E3").Value = '  -3.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.60%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3895'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3845'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.003'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.347'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.31%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '49.12'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08442'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.113'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001280'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.454'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.49%  '
$ws.Range("D17").Value = '1.638.55'
$ws.Range("E17").Value = '  -5.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06944'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.912'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.39%  '
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.73%  '
$ws.Range("D24").Value = '24.060.50'
$ws.Range("E24").Value = '  -3.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.328'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.684'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '158.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.660'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '141.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.235'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -14.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.457'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.55%  '
$ws.Range("D33").Value = '1.824.98'
$ws.Range("E33").Value = '  -5.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.131'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07995'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02904'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9603'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2686'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09210'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.463'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.940'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7592'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6891'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.476'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.088'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08366'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.255'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.41%  '
